$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 311 (shifts old rows 311-387 down to 314-390)
$ws.Rows(311).Insert()
$ws.Rows(311).Insert()
$ws.Rows(311).Insert()

# Populate the 3 newly inserted rows (311-313) with new weekly data
# Row 311
$ws.Cells.Item(311, 1).Value = 4
$ws.Cells.Item(311, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(311, 3).Value = "Los Lagos"
$ws.Cells.Item(311, 4).Value = 44551
$ws.Cells.Item(311, 5).Value = 10
$ws.Cells.Item(311, 6).Value = 100112004
$ws.Cells.Item(311, 7).Value = "Cebolla"
$ws.Cells.Item(311, 8).Value = "Morada(o)"
$ws.Cells.Item(311, 9).Value = "1a (cosecha)"
$ws.Cells.Item(311, 10).Value = 250
$ws.Cells.Item(311, 11).Value = 10000
$ws.Cells.Item(311, 12).Value = 10000
$ws.Cells.Item(311, 13).Value = 10000
$ws.Cells.Item(311, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(311, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(311, 16).Value = 556
$ws.Cells.Item(311, 17).Value = 18
$ws.Cells.Item(311, 18).Value = "Hortaliza"

# Row 312
$ws.Cells.Item(312, 1).Value = 4
$ws.Cells.Item(312, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(312, 3).Value = "Los Lagos"
$ws.Cells.Item(312, 4).Value = 44551
$ws.Cells.Item(312, 5).Value = 10
$ws.Cells.Item(312, 6).Value = 100112004
$ws.Cells.Item(312, 7).Value = "Cebolla"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "1a (cosecha)"
$ws.Cells.Item(312, 10).Value = 800
$ws.Cells.Item(312, 11).Value = 6500
$ws.Cells.Item(312, 12).Value = 6500
$ws.Cells.Item(312, 13).Value = 6500
$ws.Cells.Item(312, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(312, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(312, 16).Value = 406
$ws.Cells.Item(312, 17).Value = 16
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Row 313
$ws.Cells.Item(313, 1).Value = 4
$ws.Cells.Item(313, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(313, 3).Value = "Los Lagos"
$ws.Cells.Item(313, 4).Value = 44551
$ws.Cells.Item(313, 5).Value = 10
$ws.Cells.Item(313, 6).Value = 100112004
$ws.Cells.Item(313, 7).Value = "Cebolla"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "1a (cosecha)"
$ws.Cells.Item(313, 10).Value = 600
$ws.Cells.Item(313, 11).Value = 8000
$ws.Cells.Item(313, 12).Value = 8000
$ws.Cells.Item(313, 13).Value = 8000
$ws.Cells.Item(313, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(313, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(313, 16).Value = 444
$ws.Cells.Item(313, 17).Value = 18
$ws.Cells.Item(313, 18).Value = "Hortaliza"
